# Implemented r-skript for mca baeume
#
# The "barrierefreiheit" criterion row (row 10) is removed from the
# "Gewichtung" sheet's weighting table. Because the
# "multifunktionale_nutzungsqualitaet" group now only has 3 remaining
# criteria (versickerung, oberflaechentemperatur, befahrbarkeit) instead
# of 4, their within_group_weight values are switched from the constant
# 0.25 to an even =1/3 split - mirroring how the "kreislauffaehigkeit"
# group (also 3 criteria) already expresses its weights.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gewichtung")

# Remove the entire "barrierefreiheit" criterion row; everything below
# shifts up by one row.
$ws.Rows.Item(10).Delete()

# Re-derive the within_group_weight formulas for the
# "multifunktionale_nutzungsqualitaet" group (now rows 7-9) as an even
# =1/3 split across its remaining three criteria.
$ws.Range("F7").Formula = "=1/3"
$ws.Range("F8:F9").Formula = "=1/3"

# The "kreislauffaehigkeit" group's formulas (now rows 10-12, shifted up
# from 11-13) are re-entered the same way so they regroup into shared
# formulas consistent with the new row numbers.
$ws.Range("F10").Formula = "=1/3"
$ws.Range("F11:F12").Formula = "=1/3"

# Reflect the cell that was selected when the workbook was last saved.
$ws.Range("F19").Select()
